$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 995, shifting rows 995:1089 down to 998:1092.
$ws.Range("A995:A997").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new price records (Palta - Peru imports).
$ws.Range("A995").Value = 7
$ws.Range("B995").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C995").Value = "Ñuble"
$ws.Range("D995").Value = 45166
$ws.Range("E995").Value = 16
$ws.Range("F995").Value = "Fruta"
$ws.Range("G995").Value = 100106
$ws.Range("H995").Value = "Oleaginosos"
$ws.Range("I995").Value = 100106002
$ws.Range("J995").Value = "Palta"
$ws.Range("K995").Value = "Hass"
$ws.Range("L995").Value = "Especial"
$ws.Range("M995").Value = 100
$ws.Range("N995").Value = 30000
$ws.Range("O995").Value = 30000
$ws.Range("P995").Value = 30000
$ws.Range("Q995").Value = "`$/bandeja 10 kilos"
$ws.Range("R995").Value = "Perú"
$ws.Range("S995").Value = 3000
$ws.Range("T995").Value = 10

$ws.Range("A996").Value = 7
$ws.Range("B996").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C996").Value = "Ñuble"
$ws.Range("D996").Value = 45166
$ws.Range("E996").Value = 16
$ws.Range("F996").Value = "Fruta"
$ws.Range("G996").Value = 100106
$ws.Range("H996").Value = "Oleaginosos"
$ws.Range("I996").Value = 100106002
$ws.Range("J996").Value = "Palta"
$ws.Range("K996").Value = "Hass"
$ws.Range("L996").Value = "Primera"
$ws.Range("M996").Value = 120
$ws.Range("N996").Value = 27000
$ws.Range("O996").Value = 27000
$ws.Range("P996").Value = 27000
$ws.Range("Q996").Value = "`$/bandeja 10 kilos"
$ws.Range("R996").Value = "Perú"
$ws.Range("S996").Value = 2700
$ws.Range("T996").Value = 10

$ws.Range("A997").Value = 7
$ws.Range("B997").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C997").Value = "Ñuble"
$ws.Range("D997").Value = 45166
$ws.Range("E997").Value = 16
$ws.Range("F997").Value = "Fruta"
$ws.Range("G997").Value = 100106
$ws.Range("H997").Value = "Oleaginosos"
$ws.Range("I997").Value = 100106002
$ws.Range("J997").Value = "Palta"
$ws.Range("K997").Value = "Hass"
$ws.Range("L997").Value = "Segunda"
$ws.Range("M997").Value = 150
$ws.Range("N997").Value = 24000
$ws.Range("O997").Value = 24000
$ws.Range("P997").Value = 24000
$ws.Range("Q997").Value = "`$/bandeja 10 kilos"
$ws.Range("R997").Value = "Perú"
$ws.Range("S997").Value = 2400
$ws.Range("T997").Value = 10
